$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: "intervention_type" header + per-trial intervention type.
# Index 0 corresponds to spreadsheet row 2 (row 1 is the header row).
$interventionTypes = @(
    "DRUG",             # row 2
    "DRUG",             # row 3
    "DEVICE",           # row 4
    "DRUG",             # row 5
    "DRUG",             # row 6
    "OTHER",            # row 7
    $null,              # row 8
    "DRUG",             # row 9
    "DRUG",             # row 10
    "DRUG",             # row 11
    "DRUG",             # row 12
    "DRUG",             # row 13
    "OTHER",            # row 14
    "OTHER",            # row 15
    "DRUG",             # row 16
    "DRUG",             # row 17
    "RADIATION",        # row 18
    "DRUG",             # row 19
    "BIOLOGICAL",       # row 20
    "DRUG",             # row 21
    "DRUG",             # row 22
    "PROCEDURE",        # row 23
    "DRUG",             # row 24
    "DRUG",             # row 25
    "OTHER",            # row 26
    "DRUG",             # row 27
    "PROCEDURE",        # row 28
    "OTHER",            # row 29
    "DRUG",             # row 30
    "DRUG",             # row 31
    "OTHER",            # row 32
    "OTHER",            # row 33
    "DRUG",             # row 34
    "DRUG",             # row 35
    "DRUG",             # row 36
    "RADIATION",        # row 37
    "DEVICE",           # row 38
    "OTHER",            # row 39
    $null,              # row 40
    "DIAGNOSTIC_TEST",  # row 41
    "RADIATION",        # row 42
    "DRUG",             # row 43
    "DRUG",             # row 44
    "GENETIC",          # row 45
    "DRUG",             # row 46
    "BIOLOGICAL",       # row 47
    "BEHAVIORAL",       # row 48
    "DRUG",             # row 49
    "DRUG",             # row 50
    "DRUG",             # row 51
    "DIAGNOSTIC_TEST",  # row 52
    "DRUG",             # row 53
    "DRUG",             # row 54
    "BEHAVIORAL",       # row 55
    "OTHER",            # row 56
    "OTHER",            # row 57
    "OTHER",            # row 58
    "DRUG",             # row 59
    "DRUG",             # row 60
    "OTHER",            # row 61
    "OTHER",            # row 62
    "PROCEDURE",        # row 63
    $null,              # row 64
    "DRUG",             # row 65
    "PROCEDURE",        # row 66
    "RADIATION",        # row 67
    "DEVICE",           # row 68
    $null,              # row 69
    $null               # row 70
)

# Header cell: reuse the exact look of the rest of row 1 (bold, centered,
# boxed) by copying J1's format onto the new K1 header cell.
$ws.Range("K1").Value = "intervention_type"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the data rows (row 2 .. row 70). Trials with no reported
# intervention type still get an (empty, unstyled) text cell in column K
# so the column is fully populated down to row 70.
for ($i = 0; $i -lt $interventionTypes.Count; $i++) {
    $value = $interventionTypes[$i]
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 11)
    if ($null -ne $value) {
        $cell.Value = $value
    } else {
        $cell.Formula = "'"
        $cell.Style = "Normal"
    }
}

$ws.Range("A1:K70").Select()
